$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.262.18'
$ws.Range('E2').Value = '  -1.01%  '
$ws.Range('D3').Value = '2.241.25'
$ws.Range('E3').Value = '  -1.06%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = "'246.35"
$ws.Range('E5').Value = '  -1.36%  '
$ws.Range('E6').Value = '  -1.95%  '
$ws.Range('D7').Value = "'74.48"
$ws.Range('E7').Value = '  -3.59%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('E9').Value = '  -4.09%  '
$ws.Range('D10').Value = "'41.81"
$ws.Range('E10').Value = '  +4.26%  '
$ws.Range('E11').Value = '  -2.85%  '
$ws.Range('E12').Value = '  -2.68%  '
$ws.Range('E13').Value = '  -2.78%  '
$ws.Range('D14').Value = "'14.50"
$ws.Range('D15').Value = "'0.850"
$ws.Range('E15').Value = '  -1.52%  '
$ws.Range('D16').Value = '2.256.55'
$ws.Range('E16').Value = '  -0.81%  '
$ws.Range('D17').Value = '42.063.90'
$ws.Range('E17').Value = '  -1.19%  '
$ws.Range('D18').Value = '0.0₃0985'
$ws.Range('E18').Value = '  -0.70%  '
$ws.Range('D20').Value = "'72.00"
$ws.Range('E20').Value = '  -0.12%  '
$ws.Range('D21').Value = "'2.22"
$ws.Range('E21').Value = '  +3.54%  '
$ws.Range('D22').Value = "'231.24"
$ws.Range('E22').Value = '  -1.87%  '
$ws.Range('D23').Value = "'8.84"
$ws.Range('E23').Value = '  +38.11%  '
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').Value = "'11.27"
$ws.Range('E25').Value = '  -0.31%  '
$ws.Range('D26').Value = "'3.62"
$ws.Range('E26').Value = '  -3.38%  '
$ws.Range('D27').Value = "'2.30"
$ws.Range('E27').Value = '  -2.83%  '
$ws.Range('E28').Value = '  -1.74%  '
$ws.Range('D29').Value = "'169.07"
$ws.Range('E29').Value = '  +0.91%  '
$ws.Range('D30').Value = "'20.60"
$ws.Range('E30').Value = '  -1.59%  '
$ws.Range('E31').Value = '  -4.90%  '
$ws.Range('E32').Value = '  -3.56%  '
$ws.Range('D33').Value = "'30.18"
$ws.Range('E33').Value = '  -3.54%  '
$ws.Range('E34').Value = '  -1.66%  '
$ws.Range('D35').Value = "'5.20"
$ws.Range('E35').Value = '  +10.08%  '
$ws.Range('D36').Value = "'4.46"
$ws.Range('E36').Value = '  -2.06%  '
$ws.Range('E37').Value = '  +1.22%  '
$ws.Range('D38').Value = "'13.54"
$ws.Range('E38').Value = '  -1.76%  '
$ws.Range('E40').Value = '  -1.58%  '
$ws.Range('D41').Value = "'62.22"
$ws.Range('E41').Value = '  +1.48%  '
$ws.Range('E42').Value = '  -2.39%  '
$ws.Range('D43').Value = "'106.59"
$ws.Range('E43').Value = '  -1.54%  '
$ws.Range('E44').Value = '  +1.98%  '
$ws.Range('D45').Value = "'8.63"
$ws.Range('E45').Value = '  -2.65%  '
$ws.Range('E46').Value = '  -0.24%  '
$ws.Range('E47').Value = '  -3.07%  '
$ws.Range('D48').Value = "'4.30"
$ws.Range('E48').Value = '  -6.93%  '
$ws.Range('E49').Value = '  -0.73%  '
$ws.Range('D50').Value = "'2.25"
$ws.Range('E50').Value = '  -0.02%  '
$ws.Range('E51').Value = '  +0.27%  '

$ws.Range('D5').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
